$wb = $excel.ActiveWorkbook

# --- 1) "Active" sheet: drop the msi-install task (Id 22) ----------------
# It was fixed (see updated note) and is being re-filed under "Inactive".
$active = $wb.Worksheets.Item("Active")
$active.Rows.Item(2).Delete()

# --- 2) "Inactive" sheet: re-insert it as a Done task with the full story
$inactive = $wb.Worksheets.Item("Inactive")
$inactive.Rows.Item(2).Insert()
$inactive.Range("A2:F2").Style = $inactive.Range("A3:F3").Style

$inactive.Range("A2").Value = 22
$inactive.Range("B2").Value = "get the msi installed program to run again - it won't run after installation`nUPDATE`nthe problem was that I manually set the WindowsForm icon file in the code and the file could not be found`ntook that out since I set the icon file through project properties`nruns fine now"
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = "Bug"

$inactive.Range("E2").NumberFormat = "@"
$inactive.Range("E2").Value = "3/2/2018"
$inactive.Range("E2").Style = $inactive.Range("E3").Style

$inactive.Range("F2").NumberFormat = "@"
$inactive.Range("F2").Value = "4/11/2018"
$inactive.Range("F2").Style = $inactive.Range("F3").Style

# Re-inserting a row with multi-line text can trigger Excel's row-height
# autofit; bring row 2 back in line with the rest of the (unheighted) sheet.
$inactive.Rows.Item(2).EntireRow.AutoFit()
